$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3d printed parts")

$ws.Range("E6").Value = "Jeroen"
$ws.Range("F6").Value = 4

$ws.Range("F6").Select()
